$wb = $excel.ActiveWorkbook

# --- Sheet "About": rework the India/geothermal notes at the bottom ---
$ws1 = $wb.Worksheets.Item("About")

# Insert two new rows after the existing row 33 (pushes old rows 35-38 down to 37-40)
$ws1.Rows.Item(34).Insert()
$ws1.Rows.Item(34).Insert()

# Row 33 loses its trailing period
$ws1.Cells.Item(33, 1).Value = "For India, Flag zero for hydro because we track pumped hydro separately"

# New rows 34 and 35 with the continuation of the note about geothermal/pumped hydro
$ws1.Cells.Item(34, 1).Value = "as the geothermal plant type. Because pumped hydro plants provide half (see elec/FPC)"
$ws1.Cells.Item(35, 1).Value = "the flexibility points as peaker plants, we use a value of 0.5 for geothermal here."

# --- Sheet "BPaFF-BDTPTPF": geothermal now gets half flexibility credit ---
$ws3 = $wb.Worksheets.Item("BPaFF-BDTPTPF")
$ws3.Range("B10").Value = 0.5
